$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.240.46"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.591.18"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'212.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "'18.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.67%  "
$ws.Range("D11").Value = "'0.0852"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "1.814.88"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "1.591.76"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "'0.508"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.96%  "
$ws.Range("D16").Value = "'63.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").Value = "26.234.19"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'215.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "'9.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "'2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.24%  "
$ws.Range("D25").Value = "'144.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'6.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").Value = "'15.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("E30").Value = "  -2.77%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").Value = "'3.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").Value = "1.409.85"
$ws.Range("E33").Value = "  +5.58%  "
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("E37").Value = "  -3.81%  "
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("D40").Value = "'5.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'0.976"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "1.727.45"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'60.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("D47").Value = "'86.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("E51").Value = "  -0.02%  "
